# Update the "match_on_lang" grouped-count table:
#  - add a new "percent" column (C)
#  - refresh the DOI counts in column B
#  - populate the new percent values (stored as text) in column C

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell C1, matching the look (bold / border / centered) of the
# existing header cells by copying the formatting from B1.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C1").Value = "percent"

# Refreshed counts for the two boolean groups.
$ws.Range("B2").Value = 5223
$ws.Range("B3").Value = 4771

# New percent column, stored as text (as in the source data).
$ws.Range("C2").Formula = "'52.26"
$ws.Range("C3").Formula = "'47.74"
